# order extensibility and receipt ui
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = "NTU"
$ws.Range("C2").Value = "chicken nugget"
$ws.Range("F2").Value = "NEW"
$ws.Range("G2").Value = "Cash"

# Row 3 updates
$ws.Range("A3").Value = "D-101"
$ws.Range("C3").Value = "fries, Teat"
$ws.Range("D3").Value = "fries : spicy"
$ws.Range("E3").Value = "'false"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "Cash"

# Update selection to mirror the new active range
$ws.Range("A2:G9").Select()
